# Rename sheets
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("principal")
$ws2 = $wb.Worksheets.Item("relacion")

$ws1.Name = "tabla_banco_126"
$ws2.Name = "tabla_banco_126_rel"

# Delete the "definicion" column (column C) on the principal sheet
$ws1.Columns.Item(3).Delete()
$ws1.Rows.AutoFit()
